# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Cell B11 on sheet "Rules" changes from the text "R40" to the text "1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Use a leading apostrophe so Excel stores "1" as text (not a number),
# matching the shared-string cell type in the target workbook.
$ws.Range("B11").Value = "'1"
